$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C2" = -11.9871
    "B3" = 7.023700000000003
    "C3" = -11.87529999999999
    "B4" = 8.960299999999997
    "E8" = 15.9386
    "C9" = -10.51509999999999
    "A11" = -21.7882
    "E11" = 16.57429999999999
    "A12" = -21.5851
    "B14" = 6.865899999999998
    "E14" = 17.19129999999999
    "A15" = -21.42509999999999
    "C15" = -12.4879
    "E15" = 16.23389999999999
    "E17" = 16.39410000000002
    "C19" = -12.45840000000001
    "C20" = -11.52750000000001
    "C25" = -13.2958
    "B26" = 5.473899999999998
    "E26" = 16.33730000000001
    "A27" = -21.37839999999998
    "C27" = -13.1019
    "A28" = -21.7851
    "C28" = -13.2847
    "C30" = -11.90879999999999
    "A31" = -21.6131
    "B31" = 5.186700000000001
    "A32" = -21.51139999999999
    "C32" = -13.2651
    "B35" = 8.696200000000006
    "A36" = -20.1843
    "E36" = 16.50460000000001
    "B37" = 8.794700000000001
    "A38" = -19.6737
    "B39" = 9.526800000000007
    "B40" = 8.337999999999999
    "E42" = 16.6681
    "C44" = -13.12669999999999
    "B45" = 6.565999999999998
    "A46" = -21.535
    "C47" = -11.953
    "B52" = 5.231600000000001
    "A54" = -21.60059999999999
    "A55" = -22.5957
    "A56" = -22.0268
    "B57" = 4.442699999999995
    "C58" = -12.8135
    "C62" = -14.71410000000001
    "E64" = 17.4169
    "A67" = -21.48529999999997
    "E68" = 16.7786
    "A69" = -21.72289999999997
    "A72" = -22.17860000000003
    "A73" = -20.02489999999999
    "C77" = -11.7679
    "C78" = -11.6034
    "E79" = 18.16500000000002
    "B81" = 6.833799999999999
    "A83" = -21.49239999999999
    "B83" = 5.161300000000006
    "C84" = -13.66629999999999
    "A86" = -22.082
    "C89" = -10.2972
    "E89" = 17.37450000000002
    "A91" = -21.50049999999999
    "C91" = -10.95989999999999
    "C92" = -11.4414
    "A93" = -21.25069999999999
    "C96" = -12.5879
    "A99" = -20.23069999999999
    "B100" = 5.203399999999999
    "B102" = 8.412300000000002
    "C102" = -13.44320000000001
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()
